$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the resistor value typo: "4.7kOhm" -> "1mOhm" (R3, R4 row, Value column) ---
$ws.Range("B5").Value = "1mOhm"

# --- 2. Draw a thin gray rectangle box around two new blank rows (9 & 10) below the BOM table ---

$xlEdgeLeft   = 7
$xlEdgeTop    = 8
$xlEdgeBottom = 9
$xlEdgeRight  = 10
$xlContinuous = 1
$xlThin       = 2
$xlSolid      = 1

$white      = 16777215   # FFFFFF
$colorEdge  = 11184810   # AAAAAA - used for left/right/bottom edges
$colorTop   = 10855845   # A5A5A5 - used for the top edge

$ws.Rows.Item(9).RowHeight  = 14.7
$ws.Rows.Item(10).RowHeight = 14.7

# Row 9 is the top edge of the box: top border runs across A9:F9,
# plus a left border on A9 and a right border on F9.
foreach ($addr in @("B9","C9","D9","E9")) {
  $c = $ws.Range($addr)
  $c.Borders.Item($xlEdgeTop).Color = $colorTop
  $c.Borders.Item($xlEdgeTop).Weight = $xlThin
  $c.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
  $c.Interior.Color = $white
  $c.Interior.Pattern = $xlSolid
}

$c = $ws.Range("A9")
$c.Borders.Item($xlEdgeTop).Color = $colorTop
$c.Borders.Item($xlEdgeTop).Weight = $xlThin
$c.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$c.Borders.Item($xlEdgeLeft).Color = $colorEdge
$c.Borders.Item($xlEdgeLeft).Weight = $xlThin
$c.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$c.Interior.Color = $white
$c.Interior.Pattern = $xlSolid

$c = $ws.Range("F9")
$c.Borders.Item($xlEdgeTop).Color = $colorTop
$c.Borders.Item($xlEdgeTop).Weight = $xlThin
$c.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$c.Borders.Item($xlEdgeRight).Color = $colorEdge
$c.Borders.Item($xlEdgeRight).Weight = $xlThin
$c.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$c.Interior.Color = $white
$c.Interior.Pattern = $xlSolid

# Row 10 is the bottom edge of the box: bottom border runs across A10:F10,
# plus a left border on A10 and a right border on F10.
foreach ($addr in @("B10","C10","D10","E10")) {
  $c = $ws.Range($addr)
  $c.Borders.Item($xlEdgeBottom).Color = $colorEdge
  $c.Borders.Item($xlEdgeBottom).Weight = $xlThin
  $c.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
  $c.Interior.Color = $white
  $c.Interior.Pattern = $xlSolid
}

$c = $ws.Range("A10")
$c.Borders.Item($xlEdgeBottom).Color = $colorEdge
$c.Borders.Item($xlEdgeBottom).Weight = $xlThin
$c.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$c.Borders.Item($xlEdgeLeft).Color = $colorEdge
$c.Borders.Item($xlEdgeLeft).Weight = $xlThin
$c.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$c.Interior.Color = $white
$c.Interior.Pattern = $xlSolid

$c = $ws.Range("F10")
$c.Borders.Item($xlEdgeBottom).Color = $colorEdge
$c.Borders.Item($xlEdgeBottom).Weight = $xlThin
$c.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$c.Borders.Item($xlEdgeRight).Color = $colorEdge
$c.Borders.Item($xlEdgeRight).Weight = $xlThin
$c.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$c.Interior.Color = $white
$c.Interior.Pattern = $xlSolid

Write-Host "BOM.xlsx edits applied"
